$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Cells.Item(1, 1).Value2 = "Datos actualizados a 10 de Octubre de 2020 a las 16:05"

# Country name changes caused by re-sorting the table by case count
# (countries swap rows; each row keeps its position, the label and the
#  stats are updated to reflect the new occupant of that row)
$ws.Cells.Item(83, 1).Value2 = "Birmania"
$ws.Cells.Item(84, 1).Value2 = "Corea del Sur"
$ws.Cells.Item(107, 1).Value2 = "Tayikistan"
$ws.Cells.Item(108, 1).Value2 = "Guayana Francesa"
$ws.Cells.Item(132, 1).Value2 = "Bahamas"
$ws.Cells.Item(133, 1).Value2 = "Ruanda"
$ws.Cells.Item(134, 1).Value2 = "Republica de Africa Central"

# Updated statistics (columns B:H = Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes)
# Row 4
$ws.Cells.Item(4, 2).Value2 = 7898609
$ws.Cells.Item(4, 3).Value2 = 4131
$ws.Cells.Item(4, 4).Value2 = 5065397
$ws.Cells.Item(4, 5).Value2 = 2614494
$ws.Cells.Item(4, 6).Value2 = 0
$ws.Cells.Item(4, 7).Value2 = 70
$ws.Cells.Item(4, 8).Value2 = 218718

# Row 5
$ws.Cells.Item(5, 2).Value2 = 6985462
$ws.Cells.Item(5, 3).Value2 = 8454
$ws.Cells.Item(5, 4).Value2 = 5994698
$ws.Cells.Item(5, 5).Value2 = 883226
$ws.Cells.Item(5, 6).Value2 = 0
$ws.Cells.Item(5, 7).Value2 = 88
$ws.Cells.Item(5, 8).Value2 = 107538

# Row 21
$ws.Cells.Item(21, 2).Value2 = 338944
$ws.Cells.Item(21, 3).Value2 = 405
$ws.Cells.Item(21, 4).Value2 = 324737
$ws.Cells.Item(21, 5).Value2 = 9189
$ws.Cells.Item(21, 6).Value2 = 0
$ws.Cells.Item(21, 7).Value2 = 22
$ws.Cells.Item(21, 8).Value2 = 5018

# Row 25
$ws.Cells.Item(25, 2).Value2 = 321057
$ws.Cells.Item(25, 3).Value2 = 579
$ws.Cells.Item(25, 4).Value2 = 273500
$ws.Cells.Item(25, 5).Value2 = 37868
$ws.Cells.Item(25, 6).Value2 = 0
$ws.Cells.Item(25, 7).Value2 = 2
$ws.Cells.Item(25, 8).Value2 = 9689

# Row 36
$ws.Cells.Item(36, 2).Value2 = 127778
$ws.Cells.Item(36, 3).Value2 = 178
$ws.Cells.Item(36, 4).Value2 = 124767
$ws.Cells.Item(36, 5).Value2 = 2792
$ws.Cells.Item(36, 6).Value2 = 0
$ws.Cells.Item(36, 7).Value2 = 0
$ws.Cells.Item(36, 8).Value2 = 219

# Row 57
$ws.Cells.Item(57, 2).Value2 = 74860
$ws.Cells.Item(57, 3).Value2 = 0
$ws.Cells.Item(57, 4).Value2 = 70406
$ws.Cells.Item(57, 5).Value2 = 4182
$ws.Cells.Item(57, 6).Value2 = 0
$ws.Cells.Item(57, 7).Value2 = 1
$ws.Cells.Item(57, 8).Value2 = 272

# Row 59
$ws.Cells.Item(59, 2).Value2 = 60681
$ws.Cells.Item(59, 3).Value2 = 339
$ws.Cells.Item(59, 4).Value2 = 57587
$ws.Cells.Item(59, 5).Value2 = 2593
$ws.Cells.Item(59, 6).Value2 = 0
$ws.Cells.Item(59, 7).Value2 = 3
$ws.Cells.Item(59, 8).Value2 = 501

# Row 70
$ws.Cells.Item(70, 2).Value2 = 43945
$ws.Cells.Item(70, 3).Value2 = 281
$ws.Cells.Item(70, 4).Value2 = 37240
$ws.Cells.Item(70, 5).Value2 = 6327
$ws.Cells.Item(70, 6).Value2 = 0
$ws.Cells.Item(70, 7).Value2 = 11
$ws.Cells.Item(70, 8).Value2 = 378

# Row 71
$ws.Cells.Item(71, 2).Value2 = 41752
$ws.Cells.Item(71, 3).Value2 = 233
$ws.Cells.Item(71, 4).Value2 = 39235
$ws.Cells.Item(71, 5).Value2 = 1909
$ws.Cells.Item(71, 6).Value2 = 0
$ws.Cells.Item(71, 7).Value2 = 3
$ws.Cells.Item(71, 8).Value2 = 608

# Row 72
$ws.Cells.Item(72, 2).Value2 = 41686
$ws.Cells.Item(72, 3).Value2 = 318
$ws.Cells.Item(72, 4).Value2 = 23791
$ws.Cells.Item(72, 5).Value2 = 17272
$ws.Cells.Item(72, 6).Value2 = 0
$ws.Cells.Item(72, 7).Value2 = 2
$ws.Cells.Item(72, 8).Value2 = 623

# Row 77
$ws.Cells.Item(77, 2).Value2 = 34685
$ws.Cells.Item(77, 3).Value2 = 168
$ws.Cells.Item(77, 4).Value2 = 31536
$ws.Cells.Item(77, 5).Value2 = 2387
$ws.Cells.Item(77, 6).Value2 = 0
$ws.Cells.Item(77, 7).Value2 = 1
$ws.Cells.Item(77, 8).Value2 = 762

# Row 78
$ws.Cells.Item(78, 2).Value2 = 32082
$ws.Cells.Item(78, 3).Value2 = 444
$ws.Cells.Item(78, 4).Value2 = 25987
$ws.Cells.Item(78, 5).Value2 = 5428
$ws.Cells.Item(78, 6).Value2 = 0
$ws.Cells.Item(78, 7).Value2 = 2
$ws.Cells.Item(78, 8).Value2 = 667

# Row 83
$ws.Cells.Item(83, 2).Value2 = 26064
$ws.Cells.Item(83, 3).Value2 = 2158
$ws.Cells.Item(83, 4).Value2 = 7050
$ws.Cells.Item(83, 5).Value2 = 18416
$ws.Cells.Item(83, 6).Value2 = 0
$ws.Cells.Item(83, 7).Value2 = 32
$ws.Cells.Item(83, 8).Value2 = 598

# Row 84
$ws.Cells.Item(84, 2).Value2 = 24548
$ws.Cells.Item(84, 3).Value2 = 72
$ws.Cells.Item(84, 4).Value2 = 22624
$ws.Cells.Item(84, 5).Value2 = 1494
$ws.Cells.Item(84, 6).Value2 = 0
$ws.Cells.Item(84, 7).Value2 = 2
$ws.Cells.Item(84, 8).Value2 = 430

# Row 94
$ws.Cells.Item(94, 2).Value2 = 15452
$ws.Cells.Item(94, 3).Value2 = 64
$ws.Cells.Item(94, 4).Value2 = 11863
$ws.Cells.Item(94, 5).Value2 = 3314
$ws.Cells.Item(94, 6).Value2 = 0
$ws.Cells.Item(94, 7).Value2 = 0
$ws.Cells.Item(94, 8).Value2 = 275

# Row 95
$ws.Cells.Item(95, 2).Value2 = 15415
$ws.Cells.Item(95, 3).Value2 = 76
$ws.Cells.Item(95, 4).Value2 = 14541
$ws.Cells.Item(95, 5).Value2 = 537
$ws.Cells.Item(95, 6).Value2 = 0
$ws.Cells.Item(95, 7).Value2 = 1
$ws.Cells.Item(95, 8).Value2 = 337

# Row 107
$ws.Cells.Item(107, 2).Value2 = 10180
$ws.Cells.Item(107, 3).Value2 = 43
$ws.Cells.Item(107, 4).Value2 = 9006
$ws.Cells.Item(107, 5).Value2 = 1095
$ws.Cells.Item(107, 6).Value2 = 0
$ws.Cells.Item(107, 7).Value2 = 0
$ws.Cells.Item(107, 8).Value2 = 79

# Row 108
$ws.Cells.Item(108, 2).Value2 = 10144
$ws.Cells.Item(108, 3).Value2 = 0
$ws.Cells.Item(108, 4).Value2 = 9810
$ws.Cells.Item(108, 5).Value2 = 265
$ws.Cells.Item(108, 6).Value2 = 0
$ws.Cells.Item(108, 7).Value2 = 0
$ws.Cells.Item(108, 8).Value2 = 69

# Row 132
$ws.Cells.Item(132, 2).Value2 = 4955
$ws.Cells.Item(132, 3).Value2 = 242
$ws.Cells.Item(132, 4).Value2 = 2724
$ws.Cells.Item(132, 5).Value2 = 2125
$ws.Cells.Item(132, 6).Value2 = 0
$ws.Cells.Item(132, 7).Value2 = 4
$ws.Cells.Item(132, 8).Value2 = 106

# Row 133
$ws.Cells.Item(133, 2).Value2 = 4890
$ws.Cells.Item(133, 3).Value2 = 0
$ws.Cells.Item(133, 4).Value2 = 3555
$ws.Cells.Item(133, 5).Value2 = 1305
$ws.Cells.Item(133, 6).Value2 = 0
$ws.Cells.Item(133, 7).Value2 = 0
$ws.Cells.Item(133, 8).Value2 = 30

# Row 134
$ws.Cells.Item(134, 2).Value2 = 4853
$ws.Cells.Item(134, 3).Value2 = 0
$ws.Cells.Item(134, 4).Value2 = 1914
$ws.Cells.Item(134, 5).Value2 = 2877
$ws.Cells.Item(134, 6).Value2 = 0
$ws.Cells.Item(134, 7).Value2 = 0
$ws.Cells.Item(134, 8).Value2 = 62

# Row 145
$ws.Cells.Item(145, 2).Value2 = 3460
$ws.Cells.Item(145, 3).Value2 = 87
$ws.Cells.Item(145, 4).Value2 = 2466
$ws.Cells.Item(145, 5).Value2 = 984
$ws.Cells.Item(145, 6).Value2 = 0
$ws.Cells.Item(145, 7).Value2 = 0
$ws.Cells.Item(145, 8).Value2 = 10

# Row 180
$ws.Cells.Item(180, 2).Value2 = 477
$ws.Cells.Item(180, 3).Value2 = 0
$ws.Cells.Item(180, 4).Value2 = 461
$ws.Cells.Item(180, 5).Value2 = 16
$ws.Cells.Item(180, 6).Value2 = 0
$ws.Cells.Item(180, 7).Value2 = 0
$ws.Cells.Item(180, 8).Value2 = 0
